$d = $word.ActiveDocument

$replacements = @(
    @("692×9=6228", "810×3=2430"),
    @("630×5=3150", "649×4=2596"),
    @("553×8=4424", "278×6=1668"),
    @("980×3=2940", "650×5=3250"),
    @("871×3=2613", "413×6=2478"),
    @("948×2=1896", "857×7=5999"),
    @("168×3=504", "257×2=514"),
    @("843×2=1686", "422×9=3798"),
    @("675×5=3375", "166×3=498"),
    @("962×6=5772", "171×3=513"),
    @("984×2=1968", "302×7=2114"),
    @("682×9=6138", "810×6=4860"),
    @("347×7=2429", "751×4=3004"),
    @("234×9=2106", "995×7=6965"),
    @("455×9=4095", "533×8=4264"),
    @("933×8=7464", "106×5=530"),
    @("565×5=2825", "413×8=3304"),
    @("853×5=4265", "455×5=2275"),
    @("611×5=3055", "346×5=1730"),
    @("782×7=5474", "767×3=2301"),
    @("920×6=5520", "492×7=3444"),
    @("601×8=4808", "227×4=908"),
    @("791×2=1582", "610×6=3660"),
    @("844×4=3376", "906×4=3624"),
    @("699×2=1398", "905×3=2715")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
